$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 1.63
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 7.4
$ws.Range("J2").Value = 4.3
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.75
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 2.58

# Row 3 updates
$ws.Range("Q3").Value = 1.59
